$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.208.73"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.56%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.894.49"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.42%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.48"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4694"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.39%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4017"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.37"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07992"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9924"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.45%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.39"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.61%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.916.78"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.843"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.024"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.81"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.33%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06599"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.74%  "

$ws.Range("E19").Value = "  +0.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.43"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.229.26"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.57%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.490"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.56"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.30%  "

$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.100.08"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.81%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.69"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.62"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.104"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +11.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.081"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "116.88"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.048"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09417"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.11%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.391"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.19%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.541"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.66%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.337"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.99%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06060"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02235"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.174"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.046"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.53%  "

$ws.Range("E41").Value = "  +0.58%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1820"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.09%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.469"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +7.78%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.00"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.92%  "

$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.270"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.99%  "

$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07681"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.09"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.35%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5460"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.16%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.897"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.34"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.51%  "

$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.88"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.06%  "

